# Update the watchlist-style ticker grid on Sheet1.
# Columns: A=index (unchanged), B=Buying Opportunity, C=support Zone,
#          D=long buildup, E=Short buildup, F=FII ENTERING
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "NSE:ALLCARGO"
$ws.Range("C2").Value = "NSE:ACI"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "NSE:IEX"
$ws.Range("F2").Value = "NSE:CAMS"

# Row 3
$ws.Range("B3").Value = "NSE:ANGELONE"
$ws.Range("C3").Value = "NSE:BYKE"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "NSE:MANKIND"
$ws.Range("F3").Value = "NSE:INDIANB"

# Row 4
$ws.Range("B4").Value = "NSE:AVANTIFEED"
$ws.Range("C4").Value = "NSE:CORDSCABLE"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "NSE:PAGEIND"
$ws.Range("F4").Value = "NSE:KFINTECH"

# Row 5
$ws.Range("B5").Value = "NSE:CAMPUS"
$ws.Range("C5").Value = "NSE:DIACABS"
$ws.Range("D5").Value = ""

# Row 6 (B6 "NSE:CYIENTDLM" stays as-is)
$ws.Range("C6").Value = "NSE:EIMCOELECO"
$ws.Range("D6").Value = ""

# Row 7
$ws.Range("B7").Value = "NSE:INDIANB"
$ws.Range("C7").Value = "NSE:FINPIPE"
$ws.Range("D7").Value = ""

# Row 8
$ws.Range("B8").Value = "NSE:JUSTDIAL"
$ws.Range("C8").Value = "NSE:GRAPHITE"
$ws.Range("D8").Value = ""

# Row 9
$ws.Range("B9").Value = "NSE:KAJARIACER"
$ws.Range("C9").Value = "NSE:HBSL"
$ws.Range("D9").Value = ""

# Row 10
$ws.Range("B10").Value = "NSE:KFINTECH"
$ws.Range("C10").Value = "NSE:MOL"

# Row 11
$ws.Range("B11").Value = "NSE:KRISHANA"
$ws.Range("C11").Value = "NSE:NITIRAJ"

# Row 12
$ws.Range("B12").Value = "NSE:LICNETFN50"
$ws.Range("C12").Value = "NSE:NRAIL"

# Row 13
$ws.Range("B13").Value = "NSE:NITINSPIN"
$ws.Range("C13").Value = "NSE:PILANIINVS"

# Row 14
$ws.Range("B14").Value = "NSE:PRIVISCL"
$ws.Range("C14").Value = "NSE:ROSSELLIND"

# Row 15
$ws.Range("C15").Value = "NSE:SAKHTISUG"

# The former row 16 (index 14, NSE:RUCHIRA) is removed entirely,
# shrinking the used range from A1:F16 down to A1:F15.
$ws.Rows(16).Delete()
